$wb = $excel.ActiveWorkbook

# --- Sheet 1: ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H33").Value = 1401
$ws.Range("I33").Value = 1340.8462
$ws.Range("K33").Value = 1340.8462
$ws.Range("M33").Value = -1111.8462
$ws.Range("H103").Value = 91288.82000000001
$ws.Range("I103").Value = 166913
$ws.Range("J103").Value = 539.8
$ws.Range("K103").Value = 500739
$ws.Range("L103").Value = 1619.4
$ws.Range("M103").Value = -500153
$ws.Range("N103").Value = -2791.4
$ws.Range("H116").Value = 2878.3044
$ws.Range("I116").Value = 2740.3333
$ws.Range("J116").Value = 3137
$ws.Range("K116").Value = 2740.3333
$ws.Range("L116").Value = 3137
$ws.Range("M116").Value = 701.6667000000002
$ws.Range("N116").Value = -10021
$ws.Range("H137").Value = 2856.889
$ws.Range("I137").Value = 3119.0417
$ws.Range("J137").Value = 759.6667
$ws.Range("K137").Value = 9357.125100000001
$ws.Range("L137").Value = 2279.0001
$ws.Range("M137").Value = -6807.125100000001
$ws.Range("N137").Value = -7379.0001
$ws.Range("H138").Value = 149193.08
$ws.Range("I138").Value = 2600.75
$ws.Range("J138").Value = 189632.34
$ws.Range("K138").Value = 7802.25
$ws.Range("L138").Value = 568897.02
$ws.Range("M138").Value = -2662.25
$ws.Range("N138").Value = -579177.02

# --- Sheet 2: ARM ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 534905.9399999999
$ws.Range("I32").Value = 682515.4399999999
$ws.Range("K32").Value = 682515.4399999999
$ws.Range("M32").Value = -682228.4399999999
$ws.Range("H61").Value = 2849.0715
$ws.Range("I61").Value = 1958
$ws.Range("J61").Value = 3517.375
$ws.Range("K61").Value = 1958
$ws.Range("L61").Value = 3517.375
$ws.Range("M61").Value = -1746
$ws.Range("N61").Value = -3941.375
$ws.Range("H132").Value = 11146.923
$ws.Range("I132").Value = 14833.667
$ws.Range("K132").Value = 44501.001
$ws.Range("M132").Value = -41971.001
$ws.Range("H136").Value = 2849.0715
$ws.Range("I136").Value = 1958
$ws.Range("J136").Value = 3517.375
$ws.Range("K136").Value = 5874
$ws.Range("L136").Value = 10552.125
$ws.Range("M136").Value = -3324
$ws.Range("N136").Value = -15652.125

# --- Sheet 3: BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H108").Value = 72821
$ws.Range("J108").Value = 72821
$ws.Range("L108").Value = 72821
$ws.Range("N108").Value = -80501
$ws.Range("H134").Value = 3912.1428
$ws.Range("I134").Value = 4128.3335
$ws.Range("J134").Value = 3750
$ws.Range("K134").Value = 12385.0005
$ws.Range("L134").Value = 11250
$ws.Range("M134").Value = -9850.000499999998
$ws.Range("N134").Value = -16320

# --- Sheet 4: CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 2426.8147
$ws.Range("I31").Value = 1086.8182
$ws.Range("J31").Value = 8322.799999999999
$ws.Range("K31").Value = 1086.8182
$ws.Range("L31").Value = 8322.799999999999
$ws.Range("M31").Value = -791.8181999999999
$ws.Range("N31").Value = -8912.799999999999
$ws.Range("H34").Value = 2426.8147
$ws.Range("I34").Value = 1086.8182
$ws.Range("J34").Value = 8322.799999999999
$ws.Range("K34").Value = 1086.8182
$ws.Range("L34").Value = 8322.799999999999
$ws.Range("M34").Value = -884.8181999999999
$ws.Range("N34").Value = -8726.799999999999
$ws.Range("H58").Value = 1633.4546
$ws.Range("I58").Value = 1533.3334
$ws.Range("J58").Value = 1671
$ws.Range("K58").Value = 1533.3334
$ws.Range("L58").Value = 1671
$ws.Range("M58").Value = -1330.3334
$ws.Range("N58").Value = -2077
$ws.Range("H136").Value = 1633.4546
$ws.Range("I136").Value = 1533.3334
$ws.Range("J136").Value = 1671
$ws.Range("K136").Value = 4600.0002
$ws.Range("L136").Value = 5013
$ws.Range("M136").Value = -2050.0002
$ws.Range("N136").Value = -10113

# --- Sheet 5: CUL ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H5").Value = 832.8333
$ws.Range("I5").Value = 669.4
$ws.Range("K5").Value = 2008.2
$ws.Range("M5").Value = -1896.2
$ws.Range("H21").Value = 2175.6667
$ws.Range("I21").Value = 1388.75
$ws.Range("J21").Value = 3749.5
$ws.Range("K21").Value = 4166.25
$ws.Range("L21").Value = 11248.5
$ws.Range("M21").Value = -3993.25
$ws.Range("N21").Value = -11594.5
$ws.Range("H40").Value = 111.9
$ws.Range("J40").Value = 117.5
$ws.Range("L40").Value = 470
$ws.Range("N40").Value = -608
$ws.Range("H49").Value = 6812.5713
$ws.Range("J49").Value = 6812.5713
$ws.Range("L49").Value = 20437.7139
$ws.Range("N49").Value = -20749.7139
$ws.Range("H86").Value = 1319.25
$ws.Range("I86").Value = 998.5
$ws.Range("K86").Value = 2995.5
$ws.Range("M86").Value = -1809.5
$ws.Range("H89").Value = 1319.25
$ws.Range("I89").Value = 998.5
$ws.Range("K89").Value = 8986.5
$ws.Range("M89").Value = -3058.5
$ws.Range("H121").Value = 1123.1372
$ws.Range("I121").Value = 746.3333
$ws.Range("J121").Value = 1146.6875
$ws.Range("K121").Value = 2238.9999
$ws.Range("L121").Value = 3440.0625
$ws.Range("M121").Value = -928.9998999999998
$ws.Range("N121").Value = -6060.0625
$ws.Range("H122").Value = 7644
$ws.Range("I122").Value = 401.9
$ws.Range("J122").Value = 25749.25
$ws.Range("K122").Value = 3617.1
$ws.Range("L122").Value = 231743.25
$ws.Range("M122").Value = -1167.1
$ws.Range("N122").Value = -236643.25
$ws.Range("H131").Value = 844.53845
$ws.Range("J131").Value = 1063.579
$ws.Range("L131").Value = 3190.737
$ws.Range("N131").Value = -13270.737
$ws.Range("H135").Value = 832.8333
$ws.Range("I135").Value = 669.4
$ws.Range("K135").Value = 6024.599999999999
$ws.Range("M135").Value = -3489.599999999999

# --- Sheet 6: GSM ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H107").Value = 614.9
$ws.Range("I107").Value = 675
$ws.Range("J107").Value = 580.1053000000001
$ws.Range("K107").Value = 675
$ws.Range("L107").Value = 580.1053000000001
$ws.Range("M107").Value = 1245
$ws.Range("N107").Value = -4420.1053
$ws.Range("H122").Value = 4860.4443
$ws.Range("I122").Value = 3650.2
$ws.Range("K122").Value = 10950.6
$ws.Range("M122").Value = -8500.599999999999

# --- Sheet 7: LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H18").Value = 63004
$ws.Range("J18").Value = 63004
$ws.Range("L18").Value = 63004
$ws.Range("N18").Value = -63348
$ws.Range("H19").Value = 5440
$ws.Range("I19").Value = 880
$ws.Range("J19").Value = 10000
$ws.Range("K19").Value = 880
$ws.Range("L19").Value = 10000
$ws.Range("M19").Value = -710
$ws.Range("N19").Value = -10340
$ws.Range("H24").Value = 29335.666
$ws.Range("I24").Value = 9000
$ws.Range("K24").Value = 9000
$ws.Range("M24").Value = -8657
$ws.Range("H136").Value = 11113281
$ws.Range("I136").Value = 3220
$ws.Range("J136").Value = 16668312
$ws.Range("K136").Value = 9660
$ws.Range("L136").Value = 50004936
$ws.Range("M136").Value = -7110
$ws.Range("N136").Value = -50010036

# --- Sheet 8: WVR ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H28").Value = 57679.332
$ws.Range("J28").Value = 57679.332
$ws.Range("L28").Value = 57679.332
$ws.Range("N28").Value = -58375.332
